{"js": "// Insert a new \"Subtitle\" styled paragraph right after the document's\n// title paragraph (\"Z\u00e1pis ADB\"), containing the subtitle text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\n\nconst subtitleText =\n  \"N\u00e1vod k z\u00e1pisu Archeologick\u00fdch dokumenta\u010dn\u00edch bod\u016f (ADB) v r\u00e1mci Pra\u017esk\u00e9 pam\u00e1tkov\u00e9 rezervace\";\n\nconst newParagraph = titleParagraph.insertParagraph(subtitleText, Word.InsertLocation.after);\nnewParagraph.styleBuiltIn = Word.BuiltInStyleName.subtitle;\n\nawait context.sync();\n", "ps1": "# Insert a new \"Subtitle\" styled paragraph right after the document's\n# title paragraph (\"Z\u00e1pis ADB\"), containing the subtitle text.\n$d = $word.ActiveDocument\n\n$titleRange = $d.Paragraphs(1).Range\n$titleRange.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs(2)\n$newRange = $newPara.Range\n$newRange.Text = \"N\u00e1vod k z\u00e1pisu Archeologick\u00fdch dokumenta\u010dn\u00edch bod\u016f (ADB) v r\u00e1mci Pra\u017esk\u00e9 pam\u00e1tkov\u00e9 rezervace\"\n$newPara.Style = \"Subtitle\"\n"}
